$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Label" header in H1, matching the bold/bordered header style
# used by the other header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Populate the new "Label" column (H) for the data rows, per the refit results
$labelValues = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labelValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labelValues[$i]
}

# Refreshed precision on a handful of prediction/error values from the refit
$ws.Range("D3").Value = 0.4320950645003666
$ws.Range("E3").Value = 0.4320950645003666

$ws.Range("D7").Value = 0.5148655513084809
$ws.Range("E7").Value = 0.4851344486915191

$ws.Range("D10").Value = 0.5410059939617095
$ws.Range("E10").Value = 0.4589940060382905
